# Update the pl_mw results table (Sheet1) with new values for the
# "case with 380 kV" run. Rows 2-25 correspond to data rows 0-23;
# columns B, C, D, F, G, I, J, K, L are updated. Columns A, E, H, M, N, O
# are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","F","G","I","J","K","L")
$data = @(
    @(1.184546363813467, 0.05898662757871875, 0.008625172303449613, 4.480739842585137, 0.002642504582053729, 2.901576537038522, 0.1947821690312814, 1.018797719061894, 0.3840735125701542),
    @(1.166632628660665, 0.05376422504723166, 0.008862043121221186, 4.43978460057734, 0.002646895375753631, 2.876685712847703, 0.1946319986480241, 0.995064199117337, 0.3820585161539753),
    @(1.156455430480605, 0.05059025623481261, 0.009018774222908643, 4.415878706361369, 0.002649734368667098, 2.862081270517209, 0.1945834638562687, 0.9812202663420351, 0.3810170452359571),
    @(1.152515051120702, 0.04930495615913344, 0.009085494432505392, 4.40644854856356, 0.002650927363018538, 2.856300220834896, 0.1945746787132059, 0.9757620077615172, 0.3806419254748903),
    @(1.151873258209321, 0.04909202100454024, 0.009096745810988294, 4.404901494090538, 0.002651127641500445, 2.855350564695684, 0.1945738842616116, 0.9748667382322367, 0.3805826159880965),
    @(1.156401451053938, 0.0505728894659967, 0.009019662475668433, 4.415750266382673, 0.002649750311621625, 2.862002615767636, 0.1945833008497928, 0.9811459122787483, 0.3810117865849705),
    @(1.178199250556787, 0.05717912467443398, 0.008704508801150368, 4.466360737905589, 0.002643988912743295, 2.892853160970489, 0.1947213332992952, 1.010463213452795, 0.3833381527092996),
    @(1.227460380194771, 0.07039704093600108, 0.008175573655033519, 4.57547537887325, 0.002633820371982996, 2.958752253539558, 0.19533812827612, 1.073737646487928, 0.3894511689837685),
    @(1.267623797288479, 0.08027573552971035, 0.007840582769856885, 4.661700528148572, 0.002627030699413893, 3.010490915921068, 0.1960019560033146, 1.123762224610459, 0.3948862516103731),
    @(1.286757899876562, 0.08480777784922111, 0.00769968555272893, 4.702252221140157, 0.002624088220302689, 3.034756506402132, 0.1963496387980896, 1.147290713579793, 0.3975635026628765),
    @(1.29412757892058, 0.08652954230720411, 0.007647972538795678, 4.717799554671842, 0.002622994879825753, 3.044050559153504, 0.1964878616292793, 1.156311470971588, 0.398606717802906),
    @(1.292534873923131, 0.0861584799980335, 0.00765903701721049, 4.714442644172834, 0.002623229421889283, 3.042044234804806, 0.1964578011292417, 1.154363751185116, 0.3983807360789626),
    @(1.287361723409305, 0.08494931613398649, 0.007695398252195451, 4.703527473613462, 0.002623997852053027, 3.035519023554627, 0.1963608789956979, 1.148030631449473, 0.3976487397098794),
    @(1.284209164180169, 0.08420939748063461, 0.007717884037857781, 4.696866539807871, 0.00262447125807856, 3.031535852710789, 0.196302365800527, 1.144165872426129, 0.3972041977913392),
    @(1.266390701079303, 0.07998033165776519, 0.007850021052733247, 4.659077125478433, 0.002627225929538696, 3.008919806760048, 0.1959801530007823, 1.122240118810339, 0.394715403638088),
    @(1.2556807068068, 0.07739577716918689, 0.007934018431612122, 4.636234818258458, 0.002628953193429043, 2.995232636372691, 0.1957941851552221, 1.108987124802411, 0.3932410266666579),
    @(1.249601893596832, 0.07591279758740654, 0.007983413804985018, 4.623221447857702, 0.002629960435511181, 2.98742877510071, 0.1956915238142471, 1.101437021867156, 0.3924122842582847),
    @(1.247557679212662, 0.07541129969591509, 0.00800032444691845, 4.618836787249762, 0.002630303837910625, 2.984798300786636, 0.1956575036836483, 1.098893167539529, 0.3921349993364061),
    @(1.256812391813213, 0.07767053544037594, 0.007924964831978576, 4.638653489194127, 0.002628767899473416, 2.996682553123136, 0.1958135365311549, 1.110390407682132, 0.3933959815674797),
    @(1.288877838836129, 0.0853043247622054, 0.007684673605533154, 4.706728329066152, 0.002623771578809997, 3.037432778673505, 0.1963891693338411, 1.149887808821831, 0.3978629473368329),
    @(1.310557212766611, 0.09032600445041794, 0.007537193587794144, 4.752334333000363, 0.002620628037229848, 3.064678740514807, 0.1968036280729351, 1.176348817174329, 0.4009537000901986),
    @(1.298920444950681, 0.08764283089016089, 0.007615034860197101, 4.727891378521264, 0.002622294691369617, 3.050080834382328, 0.1965789268483533, 1.162166857627653, 0.3992884471638973),
    @(1.256300512476002, 0.07754630815277608, 0.007929054524646428, 4.637559637487129, 0.002628851626721793, 2.996026843074105, 0.1958047745225571, 1.10975576788374, 0.393325867593461),
    @(1.213436511276228, 0.0667923467791951, 0.008309199585091243, 4.544896146713768, 0.002636451078849207, 2.94034348232168, 0.1951342499620772, 1.055999967789177, 0.3876315695066097)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $col = $cols[$j]
        $ws.Range("$col$row").Value = $rowValues[$j]
    }
}

Write-Output "Updated $($data.Length) rows across $($cols.Length) columns."
